# Swap the contents of paired rows (columns B through AD only; column A,
# the sequential row id, stays in place) to reflect the re-ordering of
# match records described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param(
        [int]$RowA,
        [int]$RowB
    )

    $rangeA = $ws.Range("B$RowA`:AD$RowA")
    $rangeB = $ws.Range("B$RowB`:AD$RowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}

# NOTE: this engine does not bind named (-Name value) parameters, so call
# the helper positionally.
Swap-Rows 36 37
Swap-Rows 99 100
Swap-Rows 189 190
